# Append the new daily resale-number snapshot row (row 27) to the
# "CityResaleNum" sheet, matching the 2023-06-07 16:14 update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and D hold values that look numeric/date-like ("2023-06-07",
# "23") but must stay plain text, same as every other row in this sheet.
# Force text formatting before assigning so Excel doesn't silently convert
# them to a date serial number / a real number.
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "2023-06-07"
$ws.Range("A27").NumberFormat = "General"

$ws.Range("B27").Value = "16:10:30"
$ws.Range("C27").Value = "Wednesday"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23"
$ws.Range("D27").NumberFormat = "General"

# Remaining columns are plain numeric city figures.
$ws.Range("E27").Value = 120704
$ws.Range("F27").Value = 134272
$ws.Range("G27").Value = 159942
$ws.Range("H27").Value = 130580
$ws.Range("I27").Value = 175382
$ws.Range("J27").Value = 112704
$ws.Range("K27").Value = 200670
$ws.Range("L27").Value = 220558
$ws.Range("M27").Value = 172328
$ws.Range("N27").Value = 119779
$ws.Range("O27").Value = 38504
$ws.Range("P27").Value = 34534
$ws.Range("Q27").Value = 50643
$ws.Range("R27").Value = -1
$ws.Range("S27").Value = 36695
$ws.Range("T27").Value = -1
